$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1 (col J & K): same bold "LAT_model" header used in C1/H1/I1
$ws.Range("J1").Value = "LAT_model"
$ws.Range("K1").Value = "LAT_model"

# Header row 2 (col J & K): new model-name labels (order matters for shared-string ids)
$ws.Range("J2").Value = "逐层对抗训练模型(fgsm.eps = 0.3)"
$ws.Range("K2").Value = "逐层对抗训练模型(change the bp method and fp method)"

# Make the brand-new K1/K2 cells bold (with the same font family/charset as the
# rest of the bold header cells) so they share the existing header style
$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").Font.Charset = 134
$ws.Range("K1").Font.Family = 3

$ws.Range("K2").Font.Bold = $true
$ws.Range("K2").Font.Charset = 134
$ws.Range("K2").Font.Family = 3

# New result columns (rows 5-8) for the two new training-method results
$ws.Range("J5").Value = 0.96
$ws.Range("K5").Value = 0.99

$ws.Range("J6").Value = 0.92
$ws.Range("K6").Value = 0.97

$ws.Range("J7").Value = 0.84
$ws.Range("K7").Value = 0.89

$ws.Range("J8").Value = 0.69
$ws.Range("K8").Value = 0.59

# Match the author's final viewport/selection state
$ws.Range("K8").Select()
$excel.ActiveWindow.ScrollColumn = 9
